# Updated cryptos list (price/volume refresh), generated from commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that looks numeric (e.g. "1.00", "65.50",
# "72.362.69" with grouping dots). Prefix with an apostrophe so Excel
# keeps it as literal text instead of normalizing it into a number.

$ws.Range('D2').Value = '''72.362.69'
$ws.Range('E2').Value = '  +0.76%  '
$ws.Range('D3').Value = '''4.022.45'
$ws.Range('E3').Value = '  -0.29%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '''527.04'
$ws.Range('E5').Value = '  +1.71%  '
$ws.Range('D6').Value = '''150.58'
$ws.Range('E6').Value = '  +2.24%  '
$ws.Range('D7').Value = '''0.701'
$ws.Range('E7').Value = '  +13.26%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '''0.749'
$ws.Range('E9').Value = '  +2.05%  '
$ws.Range('E10').Value = '  -1.35%  '
$ws.Range('B11').Value = 'ShibaInu'
$ws.Range('C11').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D11').Value = '''0.0000323'
$ws.Range('E11').Value = '  -2.96%  '
$ws.Range('B12').Value = 'Avalanche'
$ws.Range('C12').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D12').Value = '''50.15'
$ws.Range('E12').Value = '  +6.97%  '
$ws.Range('D13').Value = '''10.73'
$ws.Range('E13').Value = '  -0.22%  '
$ws.Range('D14').Value = '''4.667.03'
$ws.Range('E14').Value = '  -0.23%  '
$ws.Range('D15').Value = '''3.998.87'
$ws.Range('E15').Value = '  -0.96%  '
$ws.Range('D16').Value = '''14.04'
$ws.Range('E16').Value = '  -1.33%  '
$ws.Range('E17').Value = '  -2.58%  '
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').Value = '''1.18'
$ws.Range('E19').Value = '  -2.35%  '
$ws.Range('D20').Value = '''72.419.06'
$ws.Range('E20').Value = '  +0.91%  '
$ws.Range('D21').Value = '''428.41'
$ws.Range('E21').Value = '  -2.08%  '
$ws.Range('D22').Value = '''97.42'
$ws.Range('E22').Value = '  +1.03%  '
$ws.Range('D23').Value = '''3.49'
$ws.Range('E23').Value = '  -0.54%  '
$ws.Range('E24').Value = '  +3.54%  '
$ws.Range('D25').Value = '''14.32'
$ws.Range('E25').Value = '  -1.59%  '
$ws.Range('D26').Value = '''11.13'
$ws.Range('E26').Value = '  -6.91%  '
$ws.Range('D27').Value = '''10.87'
$ws.Range('E27').Value = '  -2.97%  '
$ws.Range('D28').Value = '''3.70'
$ws.Range('E28').Value = '  +20.29%  '
$ws.Range('D29').Value = '''5.88'
$ws.Range('E29').Value = '  +1.81%  '
$ws.Range('D30').Value = '''36.69'
$ws.Range('E30').Value = '  -0.43%  '
$ws.Range('D31').Value = '''7.38'
$ws.Range('E31').Value = '  +2.12%  '
$ws.Range('E32').Value = '  +1.93%  '
$ws.Range('D33').Value = '''13.42'
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('D34').Value = '''683.37'
$ws.Range('E34').Value = '  -2.44%  '
$ws.Range('D35').Value = '''48.56'
$ws.Range('E35').Value = '  +20.13%  '
$ws.Range('D36').Value = '''65.50'
$ws.Range('E36').Value = '  -3.96%  '
$ws.Range('D37').Value = '''0.445'
$ws.Range('E37').Value = '  +1.95%  '
$ws.Range('E38').Value = '  -1.66%  '
$ws.Range('D39').Value = '''0.0₃0824'
$ws.Range('E39').Value = '  -8.31%  '
$ws.Range('B40').Value = 'ThetaToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D40').Value = '''3.41'
$ws.Range('E40').Value = '  -10.05%  '
$ws.Range('B41').Value = 'WEMIXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').Value = '''3.41'
$ws.Range('E41').Value = '  +7.97%  '
$ws.Range('E42').Value = '  +0.16%  '
$ws.Range('D43').Value = '''1.00'
$ws.Range('E43').Value = '  +0.32%  '
$ws.Range('E44').Value = '  +0.86%  '
$ws.Range('D45').Value = '''0.150'
$ws.Range('E45').Value = '  +4.01%  '
$ws.Range('D46').Value = '''9.83'
$ws.Range('E46').Value = '  +9.26%  '
$ws.Range('D47').Value = '''3.44'
$ws.Range('E47').Value = '  -2.23%  '
$ws.Range('E48').Value = '  -4.67%  '
$ws.Range('D49').Value = '''3.02'
$ws.Range('E49').Value = '  -4.07%  '
$ws.Range('B50').Value = 'FLOKI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D50').Value = '''0.000267'
$ws.Range('E50').Value = '  -1.97%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').Value = '''143.83'
$ws.Range('E51').Value = '  +0.96%  '
